$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("H4").Value = 3.6
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8

# Row 5 updates
$ws.Range("G5").Value = 3.1
$ws.Range("H5").Value = 3.7
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 3.55
$ws.Range("K5").Value = 2.25
$ws.Range("L5").Value = 2.55
$ws.Range("AA5").Value = 1.6
$ws.Range("AB5").Value = 2.07
$ws.Range("AC5").Value = 11.5
$ws.Range("AD5").Value = 17.5
$ws.Range("AE5").Value = 11.25
$ws.Range("AF5").Value = 40
$ws.Range("AG5").Value = 25
$ws.Range("AH5").Value = 30
$ws.Range("AJ5").Value = 7.3
$ws.Range("AK5").Value = 13.5
$ws.Range("AM5").Value = 350
$ws.Range("AN5").Value = 9
$ws.Range("AO5").Value = 10.5
$ws.Range("AQ5").Value = 18
$ws.Range("AR5").Value = 15

$wb.Save()
